$wb = $excel.ActiveWorkbook

# Rename the "aquisicoes" worksheet to "Folha1"
$sheet = $wb.Worksheets.Item("aquisicoes")
$sheet.Name = "Folha1"

# Move the selection/active cell on that sheet to G15
$sheet.Activate()
$sheet.Range("G15").Select()
